$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 130, pushing the existing row 130 (and below) down to 131.
$ws.Rows.Item(130).Insert()

# New row 130 data (weekly update for Fruta / hortaliza)
$ws.Cells.Item(130, 1).Value = 11
$ws.Cells.Item(130, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(130, 3).Value = "Bíobío"
$ws.Cells.Item(130, 4).Value = 44939
$ws.Cells.Item(130, 4).NumberFormat = $ws.Cells.Item(131, 4).NumberFormat
$ws.Cells.Item(130, 5).Value = 8
$ws.Cells.Item(130, 6).Value = "Fruta"
$ws.Cells.Item(130, 7).Value = 100101
$ws.Cells.Item(130, 8).Value = "Berries"
$ws.Cells.Item(130, 9).Value = 100101001
$ws.Cells.Item(130, 10).Value = "Arándano (blue)"
$ws.Cells.Item(130, 11).Value = "Sin especificar"
$ws.Cells.Item(130, 12).Value = "Primera"
$ws.Cells.Item(130, 13).Value = 270
$ws.Cells.Item(130, 14).Value = 2000
$ws.Cells.Item(130, 15).Value = 2500
$ws.Cells.Item(130, 16).Value = 2222
$ws.Cells.Item(130, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(130, 18).Value = "Provincia de Curicó"
$ws.Cells.Item(130, 19).Value = 1111
$ws.Cells.Item(130, 20).Value = 2
